# Applies the 2025-09-18 12:35:09 JST refresh of the lancers.jp listing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove the existing hyperlink objects before rewriting the table so stale
# relationship ids are not left pointing at the wrong rows.
$ws.Hyperlinks.Delete()

# Clear out the previously-populated data rows (old table was A2:H10).
$ws.Range("A2:H10").Clear()

$timestamp = '2025-09-18 12:35:09'

$rows = @(
  @{ B = 'AIチャットボットのβ版テスト参加者募集!'; C = 'システム開発'; D = '10,000 円 ~ 20,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395790'; G = 295; H = '🔥AI,Ai' },
  @{ B = '【急募】JUSTDBとOPERAcloudのAPI連携開発者募集'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5396169'; G = 250; H = '🔥API ◆開発' },
  @{ B = '【急募】入力ミス防止のためのkintone Javascript開発者募集'; C = 'システム開発'; D = '20,000 円 ~ 50,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395671'; G = 128; H = '★Java ◆開発' },
  @{ B = '初回 Laravel Livewireを使ったWebシステム開発の募集'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395502'; G = 125; H = '◆開発,システム開発' },
  @{ B = '2026年度新入社員研修Javaサブ講師'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395710'; G = 85; H = '★Java' },
  @{ B = '2026年度新入社員研修Javaサブ講師 (4~6月)'; C = 'システム開発'; D = '1,000,000 円 ~ 3,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395713'; G = 85; H = '★Java' },
  @{ B = '【急募】iOSアプリのAdMobメディエーション入札接続とeCPM改善'; C = 'システム開発'; D = '50,000 円 ~ 100,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395931'; G = 38; H = '◇アプリ' },
  @{ B = '【急募】HP保守管理とSEO対策の専門家を探しています'; C = 'システム開発'; D = '20,000 円 ~ 50,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5396003'; G = 33; H = '◇管理' },
  @{ B = '【急募】ストアーズ連携LINE予約サイト制作のプロを探しています!'; C = 'システム開発'; D = '20,000 円 ~ 50,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395825'; G = 33; H = '◇サイト' },
  @{ B = '【継続案件|お気軽にご応募ください】WebシステムのQAエンジニア募集'; C = 'システム開発'; D = '300,000 円 ~ 500,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395741'; G = 40; H = $null },
  @{ B = '社内のFAQの構築'; C = 'システム開発'; D = '500,000 円 ~ 1,000,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5396173'; G = 25; H = $null },
  @{ B = '〖リモート可〗Delphiエンジニア募集'; C = 'システム開発'; D = '300,000 円 ~ 500,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5341051'; G = 25; H = $null },
  @{ B = '【急募】ウェブプラットフォームのMVP制作を依頼します!'; C = 'システム開発'; D = '100,000 円 ~ 200,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5396017'; G = 18; H = $null },
  @{ B = 'MT5用EA(ex5)ファイルのデコンパイル'; C = 'システム開発'; D = '20,000 円 ~ 50,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5396009'; G = 13; H = $null },
  @{ B = '【急募】JotformとGoogleスプレッドシート連携のエキスパート募集!'; C = 'システム開発'; D = '5,000 円 ~ 10,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395809'; G = 10; H = $null },
  @{ B = '【MT4】ゴールドの取引を行うEAのサンプルソース納品'; C = 'システム開発'; D = '10,000 円 ~ 20,000 円 / 固定'; E = '期限情報なし'; F = 'https://www.lancers.jp/work/detail/5395799'; G = 10; H = $null }
)

$rowNum = 2
foreach ($row in $rows) {
    $ws.Cells.Item($rowNum, 1).Value = $timestamp
    $ws.Cells.Item($rowNum, 2).Value = $row.B
    $ws.Cells.Item($rowNum, 3).Value = $row.C
    $ws.Cells.Item($rowNum, 4).Value = $row.D
    $ws.Cells.Item($rowNum, 5).Value = $row.E
    $ws.Cells.Item($rowNum, 6).Value = $row.F
    $ws.Cells.Item($rowNum, 7).Value = $row.G
    if ($row.H -ne $null) {
        $ws.Cells.Item($rowNum, 8).Value = $row.H
    }
    $linkRange = $ws.Cells.Item($rowNum, 6)
    $ws.Hyperlinks.Add($linkRange, $row.F) | Out-Null
    $rowNum = $rowNum + 1
}

$lastRow = $rowNum - 1
$ws.Range("A1").Select()
